# Apply the "Tunniste" column addition to the Sheet1 header row and
# set up the page setup (paper size / orientation) as produced by a
# normal Excel save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cell (auto-creates the shared string "Tunniste").
$ws.Range("E1").Value = "Tunniste"

# Match the formatting of the other header cells (style index 1).
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the active selection to the newly added cell.
$ws.Range("E1").Select()

# Set page setup (paper size 9 = A4, portrait orientation).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
